# Adding support for irregular nouns
# Appends 15 new paradigm rows (33-47) to the "Tabla2" table on Sheet1,
# resizes the table/autofilter, and updates the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data rows -----------------------------------------------------
# Column layout: A = paradigm, B = notes (left blank here), C = gen_sing,
# D = plural. Values are entered in the same order the source workbook's
# shared-string table was built in (row 40's paradigm/plural are filled
# before rows 38-39, then its gen_sing note is filled in afterwards).

$ws.Range("A33").Value = "γέρου"
$ws.Range("D33").Value = "γέροι"

$ws.Range("A34").Value = "αθή"
$ws.Range("D34").Value = "αθήνε"

$ws.Range("A35").Value = "ούθι"
$ws.Range("D35").Value = "ουθίουνε"

$ws.Range("A36").Value = "ψιλέ"
$ws.Range("C36").Value = "ψιού"
$ws.Range("D36").Value = "ψιλ̣οί"

$ws.Range("A37").Value = "κούε"
$ws.Range("C37").Value = "κουνέ"
$ws.Range("D37").Value = "κούν̇οι"

$ws.Range("A40").Value = "άμπελε"
$ws.Range("D40").Value = "άμπελε"

$ws.Range("A38").Value = "μάτη"
$ws.Range("C38").Value = "μάτη / ματερί"
$ws.Range("D38").Value = "ματέρε"

$ws.Range("A39").Value = "σάτη"
$ws.Range("C39").Value = "σάτη / σατερί"
$ws.Range("D39").Value = "σατέρε"

$ws.Range("C40").Value = "άμπελε / άμπελ̣ή"

$ws.Range("A41").Value = "κρόπο"
$ws.Range("D41").Value = "κρόπε"

$ws.Range("A42").Value = "κόκαλε"
$ws.Range("D42").Value = "κόκα"

$ws.Range("A43").Value = "τσ̌έρβουλε"
$ws.Range("D43").Value = "τσ̌έρβα"

$ws.Range("A44").Value = "π̇ιτόκαλε"
$ws.Range("D44").Value = "π̇ιτόκα"

$ws.Range("A45").Value = "κάλ̣ι"
$ws.Range("D45").Value = "κάβα"

$ws.Range("A46").Value = "μάλ̣ι"
$ws.Range("D46").Value = "μάβα"

$ws.Range("A47").Value = "άι"
$ws.Range("D47").Value = "άζα"

# --- Formatting: match the existing data rows (14pt font -> row ht 18.75)
$ws.Range("A33:D47").Font.Size = 14

# --- Grow the "Tabla2" table / autofilter to the new extent ------------
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D47"))

# --- View state: scroll down a bit and select D27 ----------------------
$ws.Range("D27").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
